# Applies the "Industry" column insertion described by the diff.
# A new column is inserted at column C (shifting Mutual Fund, Status,
# Jan_2026, Dec_2025, Oct_2025, MoM, QoQ one column to the right, from
# C:I to D:J), and is populated with per-row industry classifications.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C; this shifts the existing C:I columns to D:J,
# carrying over their values/formatting (mirrors the header style too).
$ws.Columns("C:C").Insert()

# Header for the new column - copy the header formatting (bold, centered,
# bordered) from the adjacent header cell, then set the text.
$ws.Range("D1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C1").Value = "Industry"

# Per-row "Industry" values (row number -> industry name), matching the
# order of stocks already present in column B (unchanged by the insert).
$industries = @{
    2  = "Consumable Fuels"
    3  = "Finance"
    4  = "Insurance"
    5  = "Metals & Minerals Trading"
    6  = "Power"
    7  = "Minerals & Mining"
    8  = "Minerals & Mining"
    9  = "Oil"
    10 = "Finance"
    11 = "Power"
    12 = "Aerospace & Defense"
    13 = "Power"
    14 = "Construction"
    15 = "Petroleum Products"
    16 = "Electrical Equipment"
    17 = "Petroleum Products"
    18 = "Aerospace & Defense"
    19 = "Gas"
    20 = "Petroleum Products"
    21 = "Banks"
    22 = "Finance"
    23 = "Power"
}

foreach ($row in $industries.Keys) {
    $ws.Cells.Item($row, 3).Value = $industries[$row]
}
